$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-02-21 Friday" "2025-02-22 Saturday"
Replace-Text "825×8=" "357×7="
Replace-Text "559×8=" "625×8="
Replace-Text "489×7=" "106×2="
Replace-Text "364×7=" "522×2="
Replace-Text "187×8=" "340×3="
Replace-Text "575×6=" "896×9="
Replace-Text "583×5=" "455×5="
Replace-Text "273×3=" "522×6="
Replace-Text "777×3=" "260×8="
Replace-Text "539×4=" "639×9="
Replace-Text "870×2=" "484×2="
Replace-Text "875×9=" "577×7="
Replace-Text "382×3=" "338×9="
Replace-Text "847×5=" "441×9="
Replace-Text "893×5=" "611×5="
Replace-Text "861×6=" "167×4="
Replace-Text "654×5=" "153×8="
Replace-Text "315×8=" "693×4="
Replace-Text "741×2=" "114×7="
Replace-Text "528×5=" "667×9="
Replace-Text "768×8=" "497×4="
Replace-Text "972×2=" "688×7="
Replace-Text "387×6=" "776×2="
Replace-Text "160×8=" "785×9="
Replace-Text "367×4=" "115×6="
